$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 44.857143
$ws.Range("I8").Value = 35.666668
$ws.Range("K8").Value = 107.000004
$ws.Range("M8").Value = 31.999996
$ws.Range("H9").Value = 75
$ws.Range("I9").Value = 75
$ws.Range("K9").Value = 75
$ws.Range("M9").Value = 94
$ws.Range("H16").Value = 6881.6665
$ws.Range("J16").Value = 6881.6665
$ws.Range("L16").Value = 6881.6665
$ws.Range("N16").Value = -7341.6665
$ws.Range("H17").Value = 10000000
$ws.Range("J17").Value = 10000000
$ws.Range("L17").Value = 30000000
$ws.Range("N17").Value = -30000336
$ws.Range("H21").Value = 1999.5
$ws.Range("I21").Value = 1999.5
$ws.Range("K21").Value = 1999.5
$ws.Range("M21").Value = -1531.5
$ws.Range("H23").Value = 1999.5
$ws.Range("I23").Value = 1999.5
$ws.Range("K23").Value = 1999.5
$ws.Range("M23").Value = -1765.5
$ws.Range("H38").Value = 183.66667
$ws.Range("J38").Value = 899
$ws.Range("L38").Value = 2697
$ws.Range("N38").Value = -3441
$ws.Range("H132").Value = 4855.222
$ws.Range("I132").Value = 4483.696
$ws.Range("J132").Value = 6991.5
$ws.Range("K132").Value = 13451.088
$ws.Range("L132").Value = 20974.5
$ws.Range("M132").Value = -10921.088
$ws.Range("N132").Value = -26034.5
$ws.Range("H137").Value = 2535
$ws.Range("I137").Value = 2050
$ws.Range("K137").Value = 6150
$ws.Range("M137").Value = -3600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19657.45
$ws.Range("I32").Value = 17502.393
$ws.Range("J32").Value = 79999
$ws.Range("K32").Value = 17502.393
$ws.Range("L32").Value = 79999
$ws.Range("M32").Value = -17215.393
$ws.Range("N32").Value = -80573
$ws.Range("H41").Value = 3566.3333
$ws.Range("I41").Value = 3750
$ws.Range("J41").Value = 3199
$ws.Range("K41").Value = 3750
$ws.Range("L41").Value = 3199
$ws.Range("M41").Value = -3336
$ws.Range("N41").Value = -4027
$ws.Range("H61").Value = 2977.8333
$ws.Range("I61").Value = 2447.5334
$ws.Range("J61").Value = 5629.3335
$ws.Range("K61").Value = 2447.5334
$ws.Range("L61").Value = 5629.3335
$ws.Range("M61").Value = -2235.5334
$ws.Range("N61").Value = -6053.3335
$ws.Range("H132").Value = 5797.4165
$ws.Range("I132").Value = 2513
$ws.Range("K132").Value = 7539
$ws.Range("M132").Value = -5009
$ws.Range("H136").Value = 2977.8333
$ws.Range("I136").Value = 2447.5334
$ws.Range("J136").Value = 5629.3335
$ws.Range("K136").Value = 7342.600199999999
$ws.Range("L136").Value = 16888.0005
$ws.Range("M136").Value = -4792.600199999999
$ws.Range("N136").Value = -21988.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 10000
$ws.Range("I20").Value = 10000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 10000
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -9753
$ws.Range("H107").Value = 1577.8572
$ws.Range("I107").Value = 1090.8334
$ws.Range("K107").Value = 1090.8334
$ws.Range("M107").Value = 829.1666
$ws.Range("H134").Value = 4253.3335
$ws.Range("I134").Value = 3565.1667
$ws.Range("K134").Value = 10695.5001
$ws.Range("M134").Value = -8160.500100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7763.2
$ws.Range("I58").Value = 5818.2856
$ws.Range("J58").Value = 12301.333
$ws.Range("K58").Value = 5818.2856
$ws.Range("L58").Value = 12301.333
$ws.Range("M58").Value = -5615.2856
$ws.Range("N58").Value = -12707.333
$ws.Range("H59").Value = 29976.125
$ws.Range("J59").Value = 34995.555
$ws.Range("L59").Value = 34995.555
$ws.Range("N59").Value = -37285.555
$ws.Range("H122").Value = 2808.7
$ws.Range("I122").Value = 2693.3684
$ws.Range("K122").Value = 8080.1052
$ws.Range("M122").Value = -5630.1052
$ws.Range("H132").Value = 1387.125
$ws.Range("I132").Value = 1387.125
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4161.375
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -1631.375
$ws.Range("H134").Value = 8277
$ws.Range("I134").Value = 7530.3335
$ws.Range("J134").Value = 14997
$ws.Range("K134").Value = 22591.0005
$ws.Range("L134").Value = 44991
$ws.Range("M134").Value = -20056.0005
$ws.Range("N134").Value = -50061
$ws.Range("H136").Value = 7763.2
$ws.Range("I136").Value = 5818.2856
$ws.Range("J136").Value = 12301.333
$ws.Range("K136").Value = 17454.8568
$ws.Range("L136").Value = 36903.999
$ws.Range("M136").Value = -14904.8568
$ws.Range("N136").Value = -42003.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2400
$ws.Range("J113").Value = 2400
$ws.Range("L113").Value = 7200
$ws.Range("N113").Value = -11540

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6300.8184
$ws.Range("I70").Value = 5755.857
$ws.Range("K70").Value = 5755.857
$ws.Range("M70").Value = -5485.857
$ws.Range("H73").Value = 6300.8184
$ws.Range("I73").Value = 5755.857
$ws.Range("K73").Value = 5755.857
$ws.Range("M73").Value = -4819.857
$ws.Range("H122").Value = 33664.332
$ws.Range("I122").Value = 34791.06
$ws.Range("K122").Value = 104373.18
$ws.Range("M122").Value = -101923.18
$ws.Range("H132").Value = 52
$ws.Range("I132").Value = 52
$ws.Range("K132").Value = 156
$ws.Range("M132").Value = 2374

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 5013
$ws.Range("I32").Value = 5013
$ws.Range("K32").Value = 5013
$ws.Range("M32").Value = -4696
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("N74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("N77").Value = 0
$ws.Range("H122").Value = 3753.25
$ws.Range("I122").Value = 3502.6667
$ws.Range("K122").Value = 10508.0001
$ws.Range("M122").Value = -8058.000100000001
$ws.Range("H136").Value = 6479.5557
$ws.Range("I136").Value = 3453
$ws.Range("J136").Value = 8900.799999999999
$ws.Range("K136").Value = 10359
$ws.Range("L136").Value = 26702.4
$ws.Range("M136").Value = -7809
$ws.Range("N136").Value = -31802.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 23833
$ws.Range("I39").Value = 23833
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 23833
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -23420
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 1500
$ws.Range("K122").Value = 4500
$ws.Range("M122").Value = -2050
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").ClearContents()
$ws.Range("N129").Value = 0
$ws.Range("H136").Value = 3649.5
$ws.Range("I136").Value = 3649.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10948.5
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -8398.5
